{"js": "// The author expanded the closing sentence of the last paragraph of the\n// procesboek: right after \"... de button voor de wishlist.\" a number of\n// new sentences are appended, describing the close-method / winnaar logic\n// (close method toggling the listing's active flag, migrating the model\n// after adding a winner field, and picking the auction winner).\n//\n// Strategy: locate the unique anchor \"wishlist.\" (it occurs exactly once,\n// at the very end of the document) and insert the new text immediately\n// after it, so the paragraph text becomes \"...wishlist.\" + the appended\n// text, with the rest of the document left untouched.\n\nconst additionalText =\n  \" De close method past de actieve status van de listing aan, dit was eigenlijk vrij simpel door te voeren.\" +\n  \" Ik merkte dat ik er op dit moment wel een beetje de vaart in kreeg en dat ik steeds beter begrijp hoe alle bestanden in relatie staan met elkaar.\" +\n  \" Toch was ik een half uur bezig met het zoeken naar een bepaalde error die ik kreeg, maar omdat ik mijn model had aangepast om een winner toe te voegen moest ik dit natuurlijk migraten. Dit was ik vergeten. \" +\n  \"Vervolgens kies ik de winnaar, degene die de hoogste bieding heeft gedaan, op dezelfde manier als dat ik het hoogste bod zocht bij de functie voor het plaatsen van een bod. \";\n\nconst body = context.document.body;\n\nconst searchResults = body.search(\"wishlist.\", { matchCase: true, matchWholeWord: false });\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length > 0) {\n  // The anchor is expected to be unique (end of the last paragraph); use\n  // the last match just in case it is not.\n  const anchor = searchResults.items[searchResults.items.length - 1];\n  anchor.insertText(additionalText, Word.InsertLocation.after);\n} else {\n  // Fallback: the anchor text was not found as expected (e.g. different\n  // spacing) - append to the very end of the document's last paragraph.\n  const paragraphs = body.paragraphs;\n  paragraphs.load(\"items\");\n  await context.sync();\n\n  const lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n  lastParagraph.insertText(additionalText, Word.InsertLocation.end);\n}\n\nawait context.sync();\n", "ps1": "# The author expanded the closing sentence of the last paragraph of the\n# procesboek: right after \"... de button voor de wishlist.\" a number of\n# new sentences are appended, describing the close-method / winnaar logic\n# (close method toggling the listing's active flag, migrating the model\n# after adding a winner field, and picking the auction winner).\n#\n# Strategy: use Find/Replace on the unique anchor \"wishlist.\" (it occurs\n# exactly once, at the very end of the document) and replace it with\n# itself plus the newly authored text, so the run's formatting (nl-NL\n# language) is inherited by the appended text instead of creating an\n# unformatted run.\n\n$d = $word.ActiveDocument\n\n$additionalText = \" De close method past de actieve status van de listing aan, dit was eigenlijk vrij simpel door te voeren.\" + `\n    \" Ik merkte dat ik er op dit moment wel een beetje de vaart in kreeg en dat ik steeds beter begrijp hoe alle bestanden in relatie staan met elkaar.\" + `\n    \" Toch was ik een half uur bezig met het zoeken naar een bepaalde error die ik kreeg, maar omdat ik mijn model had aangepast om een winner toe te voegen moest ik dit natuurlijk migraten. Dit was ik vergeten. \" + `\n    \"Vervolgens kies ik de winnaar, degene die de hoogste bieding heeft gedaan, op dezelfde manier als dat ik het hoogste bod zocht bij de functie voor het plaatsen van een bod. \"\n\n$anchor = \"wishlist.\"\n$replacement = $anchor + $additionalText\n\n$range = $d.Content\n$found = $range.Find.Execute($anchor, $false, $false, $false, $false, $false, $true, 1, $false, $replacement, 2)\n\nif (-not $found) {\n    # Fallback: the anchor text was not found as expected (e.g. different\n    # spacing) - append to the very end of the last paragraph instead.\n    $paragraphs = $d.Paragraphs\n    $lastParagraph = $paragraphs.Item($paragraphs.Count)\n    $endRange = $lastParagraph.Range\n    $endRange.Collapse(0)\n    $endRange.InsertAfter($additionalText)\n}\n"}
